# Update data for the week through 2022-06-19 (add data for 2022-06-27 commit)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet/tab to reflect the new "through" date
$ws.Name = "Through 2022-06-19"

# Update the label for the June row
$ws.Range("A7").Value = "June (through 06-19)"

# Row 7: June monthly totals by year (B=2015 .. I=2022)
$ws.Range("C7").Value = 22
$ws.Range("D7").Value = 41
$ws.Range("F7").Value = 33
$ws.Range("G7").Value = 71
$ws.Range("H7").Value = 73
$ws.Range("I7").Value = 98

# Row 8: Year-to-date totals by year (B=2015 .. I=2022)
$ws.Range("C8").Value = 231
$ws.Range("D8").Value = 357
$ws.Range("F8").Value = 237
$ws.Range("G8").Value = 429
$ws.Range("H8").Value = 704
$ws.Range("I8").Value = 761
